$wb = $excel.ActiveWorkbook

# --- Sheet "Sezony I liga" (sheet4.xml): add round 10 (column M) stats for the 2021/2022 season block (rows 17-20) ---
$wsIliga = $wb.Worksheets.Item("Sezony I liga")
$wsIliga.Range("M17").Value = 25
$wsIliga.Range("M18").Value = 2
$wsIliga.Range("M19").Value = 5
$wsIliga.Range("M20").Value = 3

# Update the view: move the active selection to M21
$wsIliga.Activate()
[void]$wsIliga.Range("M21").Select()

# --- Sheet "Sezony Ekstra" (sheet3.xml): add round 10 (column L) stats for the 2021/2022 season block (rows 17-20) ---
$wsEkstra = $wb.Worksheets.Item("Sezony Ekstra")
$wsEkstra.Range("L17").Value = 32
$wsEkstra.Range("L18").Value = 2
$wsEkstra.Range("L19").Value = 4
$wsEkstra.Range("L20").Value = 6

# Update the view: scroll back to column A and move the active selection to L19.
# This sheet stays the active/selected tab, matching the saved workbook state.
$wsEkstra.Activate()
$wsEkstra.Application.ActiveWindow.ScrollColumn = 1
[void]$wsEkstra.Range("L19").Select()
